$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert two new blank rows to make room for the new records ---
# Insert at row 3 (shifts old rows 3,4,5 -> 4,5,6)
$ws.Rows.Item(3).Insert()
# Insert at row 6 (shifts the row currently at 6, i.e. old row 5 "fan" entry, -> 7)
$ws.Rows.Item(6).Insert()

# --- Step 2: widen column D ---
$ws.Columns.Item(4).ColumnWidth = 30

# --- Step 3: write cell values for every data row (2..8) ---
# Row 2
$ws.Range("A2").Value = '2025-11-22 12:31:51'
$ws.Range("B2").Value = '【技術者募集】家庭用消臭デバイス「Maneki Air」開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5439445'
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-11-22 12:31:51'
$ws.Range("B3").Value = '【急募】掲示板サイト(爆サイ)自動書き込みソフト開発者募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5439484'
$ws.Range("G3").Value = 93
$ws.Range("H3").Value = '◆開発 ◇サイト'

# Row 4
$ws.Range("A4").Value = '2025-11-22 12:31:51'
$ws.Range("B4").Value = '名刺/プロフィール共有アプリ開発'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5439373'
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = '◆開発 ◇アプリ'

# Row 5
$ws.Range("A5").Value = '2025-11-22 12:31:51'
$ws.Range("B5").Value = '【緊急】海外からWordPress管理画面にログインできない問題の調査と修正'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5439402'
$ws.Range("G5").Value = 50
$ws.Range("H5").Value = '◇管理 ○WordPress'

# Row 6
$ws.Range("A6").Value = '2025-11-22 12:31:51'
$ws.Range("B6").Value = '限定公開 限定公開の仕事'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5439488'
$ws.Range("G6").Value = 25

# Row 7
$ws.Range("A7").Value = '2025-11-22 12:31:51'
$ws.Range("B7").Value = '【急募】ファン応援プラットフォームの構築をお手伝いください!'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5439395'
$ws.Range("G7").Value = 18

# Row 8
$ws.Range("A8").Value = '2025-11-22 12:31:51'
$ws.Range("B8").Value = '【Stable Diffusion】参考動画に沿って約100プロンプト構築'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5432055'
$ws.Range("G8").Value = 10

# --- Step 4: rebuild hyperlinks on column F for rows 2..8 in order ---
$hn = $ws.Hyperlinks.Count()
if ($hn -gt 0) { $ws.Range("F2").Hyperlinks.Delete() }
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5439445') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5439484') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5439373') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5439402') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5439488') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5439395') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5432055') | Out-Null

Write-Output $ws.UsedRange.Address()
